# Insert a new data row at row 195 (pushing existing rows 195:228 down to 196:229)
# and populate it with the new "Vega Modelo de Temuco - Zanahoria" record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(195).Insert()

$ws.Cells.Item(195, 1).Value = 10
$ws.Cells.Item(195, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(195, 3).Value = "La Araucanía"
$ws.Cells.Item(195, 4).Value = 44522
$ws.Cells.Item(195, 5).Value = 9
$ws.Cells.Item(195, 6).Value = 100114013
$ws.Cells.Item(195, 7).Value = "Zanahoria"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 100
$ws.Cells.Item(195, 11).Value = 8000
$ws.Cells.Item(195, 12).Value = 8000
$ws.Cells.Item(195, 13).Value = 8000
$ws.Cells.Item(195, 14).Value = '$/saco 20 kilos'
$ws.Cells.Item(195, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(195, 16).Value = 400
$ws.Cells.Item(195, 17).Value = 20
$ws.Cells.Item(195, 18).Value = "Hortaliza"
